$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as text in the source data.
# Set each target cell to Text format individually (avoiding multi-area
# Union ranges, whose NumberFormat setter only affects the first area)
# before writing the updated price so the values stay text, not numbers.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D25","D26","D27","D40","D41","D42","D43","D44","D48")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "245.03"
$ws.Range("D3").Value = "25.12"
$ws.Range("D4").Value = "5.009"
$ws.Range("D5").Value = "0.05611"
$ws.Range("D6").Value = "6.574"
$ws.Range("D7").Value = "3.007"
$ws.Range("D8").Value = "0.8098"
$ws.Range("D9").Value = "0.8416"
$ws.Range("D10").Value = "0.1336"
$ws.Range("D11").Value = "0.06941"
$ws.Range("D13").Value = "0.09394"
$ws.Range("D14").Value = "0.001507"
$ws.Range("D15").Value = "0.0005977"
$ws.Range("D16").Value = "0.006148"
$ws.Range("D17").Value = "3.498"
$ws.Range("D18").Value = "2.091"
$ws.Range("D19").Value = "0.3195"
$ws.Range("D20").Value = "0.03270"
$ws.Range("D21").Value = "0.1291"
$ws.Range("D22").Value = "3.741"
$ws.Range("D23").Value = "0.04680"
$ws.Range("D25").Value = "0.001242"
$ws.Range("D26").Value = "0.004526"
$ws.Range("D27").Value = "0.00009695"
$ws.Range("D40").Value = "0.03661"
$ws.Range("D41").Value = "0.1367"
$ws.Range("D42").Value = "0.006236"
$ws.Range("D43").Value = "0.002734"
$ws.Range("D44").Value = "0.008052"
$ws.Range("D48").Value = "0.002037"
